$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the current row 121,
# pushing the existing rows 121-127 down to 122-128 (their data is
# unchanged by the shift). Insert a fresh row at 121 so the rows below
# move down and keep their own values/formatting.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new record.
$ws.Cells.Item(121, 1).Value = 3
$ws.Cells.Item(121, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44578
$ws.Cells.Item(121, 5).Value = 5
$ws.Cells.Item(121, 6).Value = 100112030
$ws.Cells.Item(121, 7).Value = "Poroto granado"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 73
$ws.Cells.Item(121, 11).Value = 28000
$ws.Cells.Item(121, 12).Value = 29000
$ws.Cells.Item(121, 13).Value = 28479
$ws.Cells.Item(121, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(121, 15).Value = "Provincia de Talca"
$ws.Cells.Item(121, 16).Value = 1139
$ws.Cells.Item(121, 17).Value = 25
$ws.Cells.Item(121, 18).Value = "Hortaliza"
